# "Luận Cung Mệnh Vô chính Diệu, Thân Mệnh đồng cung"
#
# The last two rows of the "chị" (cung Mệnh) block on Sheet1 (B51/B52, which
# both read "Cuộc đời chị an nhàn.") get re-pointed to a new verdict text,
# and a brand-new row 53 is appended that calls out the "Thân Mệnh đồng cung
# Vô Chính Diệu" case with its own consequence text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 53: Thân/Mệnh đồng cung Vô Chính Diệu + its hardship verdict.
$ws.Range("A53").Value = "Thân và Mệnh đồng cung Vô Chính Diệu"
$ws.Range("B53").Value = "Cuộc đời cực kỳ vất vả, khổ cực. Không cậy nhờ được sự giúp đỡ của người khác, tự thân lập nghiệp."

# Column A on this sheet (rows 22-52) carries a yellow highlight; match it
# for the newly appended row.
$ws.Range("A53").Interior.Color = 65535

# B51 / B52 previously shared the text "Cuộc đời chị an nhàn." - move them to
# the new verdict.
$ws.Range("B51").Value = "Độ số an nhàn gia tăng."
$ws.Range("B52").Value = "Độ số an nhàn gia tăng."

# Reflect the author's final cursor position / scroll state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C52").Select()
